# Direct Opex Function update: CO2 locations now report fuel usage in MJ
# rather than kg. This reshuffles the last few "Usage" columns (AF:AI) so
# that tkm-N1Usage / pkmUsage move ahead of tkm-SZMUsage / tkm-N2Usage
# (tkm-N3Usage in AJ is unaffected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the header row (row 1) for columns AF:AI.
$ws.Range("AF1").Value = "tkm-N1Usage"
$ws.Range("AG1").Value = "pkmUsage"
$ws.Range("AH1").Value = "tkm-SZMUsage"
$ws.Range("AI1").Value = "tkm-N2Usage"

# Re-order the matching data row (row 2) for columns AF:AI so each value
# stays with its header.
$ws.Range("AF2").Value = 7.5
$ws.Range("AG2").Value = 858
$ws.Range("AH2").Value = 414.5
$ws.Range("AI2").Value = 24.2
